$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

# Columns A (date-like) and D (numeric-looking "Week") must stay as TEXT.
# A plain string assignment gets auto-parsed by Excel into a date serial /
# number, so force text mode via NumberFormat, assign, then strip the
# resulting explicit style back off with ClearFormats so the cell keeps
# the sheet's default (unstyled) appearance.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-09"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:10:07"

$ws.Cells.Item($row, 3).Value = "Friday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "23"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 120959
$ws.Cells.Item($row, 6).Value = 134486
$ws.Cells.Item($row, 7).Value = 160802
$ws.Cells.Item($row, 8).Value = 131998
$ws.Cells.Item($row, 9).Value = 175908
$ws.Cells.Item($row, 10).Value = 113755
$ws.Cells.Item($row, 11).Value = 201753
$ws.Cells.Item($row, 12).Value = 221995
$ws.Cells.Item($row, 13).Value = 173438
$ws.Cells.Item($row, 14).Value = 118682
$ws.Cells.Item($row, 15).Value = 38714
$ws.Cells.Item($row, 16).Value = 34342
$ws.Cells.Item($row, 17).Value = 51086
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 37065
$ws.Cells.Item($row, 20).Value = -1
